$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 213; this shifts the existing rows 213:242
# down to 214:243 (so old row 242 becomes new row 243), matching the
# target dimension A1:R243.
$ws.Rows("213:213").Insert()

# Populate the newly inserted row 213 with the new weekly record.
$ws.Range("A213").Value = 7
$ws.Range("B213").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C213").Value = "Ñuble"
$ws.Range("D213").Value = 44776
$ws.Range("E213").Value = 16
$ws.Range("F213").Value = 100112003
$ws.Range("G213").Value = "Ajo"
$ws.Range("H213").Value = "Chino"
$ws.Range("I213").Value = "Primera"
$ws.Range("J213").Value = 100
$ws.Range("K213").Value = 24000
$ws.Range("L213").Value = 25000
$ws.Range("M213").Value = 24500
$ws.Range("N213").Value = "$/caja 10 kilos"
$ws.Range("O213").Value = "China"
$ws.Range("P213").Value = 2450
$ws.Range("Q213").Value = 10
$ws.Range("R213").Value = "Hortaliza"
